$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original number formatting/style for price (D) and volume (E) columns.
# These columns store values as text (e.g. "1.000", "27.930.96") which Excel
# would otherwise auto-convert to numbers. Temporarily force Text format,
# assign the new values, then restore the original style.
$rangeD = $ws.Range("D2:D51")
$rangeE = $ws.Range("E2:E51")
$styleD = $rangeD.Style
$styleE = $rangeE.Style
$rangeD.NumberFormat = "@"
$rangeE.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '28.091.79'
$ws.Cells.Item(2, 5).Value = '  -1.46%  '
$ws.Cells.Item(3, 4).Value = '1.801.91'
$ws.Cells.Item(3, 5).Value = '  -1.31%  '
$ws.Cells.Item(4, 4).Value = '0.9971'
$ws.Cells.Item(4, 5).Value = '  -0.70%  '
$ws.Cells.Item(5, 4).Value = '310.86'
$ws.Cells.Item(5, 5).Value = '  -1.79%  '
$ws.Cells.Item(6, 4).Value = '0.9975'
$ws.Cells.Item(6, 5).Value = '  -0.58%  '
$ws.Cells.Item(7, 4).Value = '0.5147'
$ws.Cells.Item(7, 5).Value = '  -0.13%  '
$ws.Cells.Item(8, 4).Value = '0.3977'
$ws.Cells.Item(8, 5).Value = '  +2.56%  '
$ws.Cells.Item(9, 4).Value = '0.07865'
$ws.Cells.Item(9, 5).Value = '  -6.43%  '
$ws.Cells.Item(10, 4).Value = '1.099'
$ws.Cells.Item(10, 5).Value = '  -2.05%  '
$ws.Cells.Item(11, 4).Value = '41.08'
$ws.Cells.Item(11, 5).Value = '  -2.17%  '
$ws.Cells.Item(12, 4).Value = '6.300'
$ws.Cells.Item(12, 5).Value = '  -2.05%  '
$ws.Cells.Item(13, 4).Value = '0.9968'
$ws.Cells.Item(13, 5).Value = '  -0.65%  '
$ws.Cells.Item(14, 4).Value = '20.39'
$ws.Cells.Item(14, 5).Value = '  -4.05%  '
$ws.Cells.Item(15, 4).Value = '7.279'
$ws.Cells.Item(15, 5).Value = '  -3.23%  '
$ws.Cells.Item(16, 4).Value = '1.785.03'
$ws.Cells.Item(16, 5).Value = '  -2.12%  '
$ws.Cells.Item(17, 4).Value = '92.08'
$ws.Cells.Item(17, 5).Value = '  -2.38%  '
$ws.Cells.Item(18, 4).Value = '0.00001080'
$ws.Cells.Item(18, 5).Value = '  -4.62%  '
$ws.Cells.Item(19, 4).Value = '0.06527'
$ws.Cells.Item(19, 5).Value = '  -1.62%  '
$ws.Cells.Item(20, 4).Value = '0.9984'
$ws.Cells.Item(20, 5).Value = '  -0.46%  '
$ws.Cells.Item(21, 4).Value = '17.20'
$ws.Cells.Item(21, 5).Value = '  -3.24%  '
$ws.Cells.Item(22, 4).Value = '5.954'
$ws.Cells.Item(22, 5).Value = '  -2.16%  '
$ws.Cells.Item(23, 4).Value = '28.155.67'
$ws.Cells.Item(23, 5).Value = '  -1.41%  '
$ws.Cells.Item(24, 4).Value = '11.06'
$ws.Cells.Item(24, 5).Value = '  -3.23%  '
$ws.Cells.Item(25, 4).Value = '2.224'
$ws.Cells.Item(25, 5).Value = '  -2.99%  '
$ws.Cells.Item(26, 4).Value = '160.20'
$ws.Cells.Item(26, 5).Value = '  +0.06%  '
$ws.Cells.Item(27, 4).Value = '20.49'
$ws.Cells.Item(27, 5).Value = '  -3.72%  '
$ws.Cells.Item(28, 4).Value = '1.996.75'
$ws.Cells.Item(28, 5).Value = '  -1.99%  '
$ws.Cells.Item(29, 4).Value = '2.385'
$ws.Cells.Item(29, 5).Value = '  -0.82%  '
$ws.Cells.Item(30, 4).Value = '127.05'
$ws.Cells.Item(30, 5).Value = '  +0.99%  '
$ws.Cells.Item(31, 4).Value = '0.1080'
$ws.Cells.Item(31, 5).Value = '  -1.44%  '
$ws.Cells.Item(32, 4).Value = '1.047'
$ws.Cells.Item(32, 5).Value = '  -4.80%  '
$ws.Cells.Item(33, 4).Value = '3.634'
$ws.Cells.Item(33, 5).Value = '  -1.15%  '
$ws.Cells.Item(34, 4).Value = '5.532'
$ws.Cells.Item(34, 5).Value = '  -3.73%  '
$ws.Cells.Item(35, 4).Value = '0.07112'
$ws.Cells.Item(35, 5).Value = '  -7.91%  '
$ws.Cells.Item(36, 4).Value = '9.027'
$ws.Cells.Item(36, 5).Value = '  +2.94%  '
$ws.Cells.Item(37, 4).Value = '0.02311'
$ws.Cells.Item(37, 5).Value = '  -3.35%  '
$ws.Cells.Item(38, 4).Value = '0.2139'
$ws.Cells.Item(38, 5).Value = '  -3.97%  '
$ws.Cells.Item(39, 4).Value = '11.61'
$ws.Cells.Item(39, 5).Value = '  +0.91%  '
$ws.Cells.Item(40, 4).Value = '5.031'
$ws.Cells.Item(40, 5).Value = '  -4.59%  '
$ws.Cells.Item(41, 4).Value = '0.6157'
$ws.Cells.Item(41, 5).Value = '  -4.15%  '
$ws.Cells.Item(42, 4).Value = '0.9984'
$ws.Cells.Item(42, 5).Value = '  -0.33%  '
$ws.Cells.Item(43, 4).Value = '1.157'
$ws.Cells.Item(43, 5).Value = '  -2.74%  '
$ws.Cells.Item(44, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(44, 4).Value = '1.316'
$ws.Cells.Item(44, 5).Value = '  -6.20%  '
$ws.Cells.Item(45, 2).Value = 'EnergySwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(45, 4).Value = '13.15'
$ws.Cells.Item(45, 5).Value = '  -3.30%  '
$ws.Cells.Item(46, 4).Value = '0.5972'
$ws.Cells.Item(46, 5).Value = '  -3.28%  '
$ws.Cells.Item(47, 4).Value = '3.722'
$ws.Cells.Item(47, 5).Value = '  -1.95%  '
$ws.Cells.Item(48, 4).Value = '125.69'
$ws.Cells.Item(48, 5).Value = '  -1.59%  '
$ws.Cells.Item(49, 4).Value = '1.211'
$ws.Cells.Item(49, 5).Value = '  +0.15%  '
$ws.Cells.Item(50, 4).Value = '1.923'
$ws.Cells.Item(50, 5).Value = '  -4.03%  '
$ws.Cells.Item(51, 4).Value = '0.06847'
$ws.Cells.Item(51, 5).Value = '  -2.19%  '

$rangeD.Style = $styleD
$rangeE.Style = $styleE
